# Update "想去人数" (F column) figures for the 2024-08-24 snapshot regeneration.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 187
$ws1.Range("F5").Value  = 3419
$ws1.Range("F13").Value = 199
$ws1.Range("F14").Value = 35
$ws1.Range("F15").Value = 87
$ws1.Range("F16").Value = 2810
$ws1.Range("F17").Value = 1121

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 187
$ws4.Range("F6").Value  = 3419
$ws4.Range("F15").Value = 199
$ws4.Range("F16").Value = 35
$ws4.Range("F17").Value = 87
$ws4.Range("F18").Value = 2810
$ws4.Range("F19").Value = 1121
